$d = $word.ActiveDocument

# Remove the trailing "Ver no Jupiter..." / copyright footer block, along with
# the blank paragraph that separated it from the "Requisitos" section, while
# keeping the paragraph containing "LOT2013: Engenharia Bioquímica I
# (Requisito fraco)" and the blank paragraph + page-break paragraph that
# follow the removed block.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Ver no Jupiter*") {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text.Trim() -eq "") {
            $startPara = $prev
        } else {
            $startPara = $p
        }
    }

    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
